$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto snapshot values.
# D-column text is forced to remain plain text (matching the source feed formatting,
# e.g. "27.682.90" / "0.4730") by temporarily switching the cell to a text number
# format before assigning the value, then restoring its original style afterwards.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.682.90"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +0.72%  "

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.877.60"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  -0.35%  "

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "331.01"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +3.13%  "

$ws.Range("E6").Value = "  -0.44%  "

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4730"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  +5.12%  "

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3972"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +2.73%  "

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "47.79"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +0.26%  "

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08055"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("E11").Value = "  +0.53%  "

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.87"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +1.63%  "

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.878.60"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -2.09%  "

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.968"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +1.23%  "

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.179"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("E16").Value = "  -0.60%  "

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "87.15"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +1.18%  "

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001043"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +0.66%  "

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06631"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.31%  "

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.23"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("E21").Value = "  -0.27%  "

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.694.39"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +0.60%  "

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.513"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("E25").Value = "  +0.77%  "

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.097.88"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -1.75%  "

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "156.35"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +3.38%  "

$ws.Range("E28").Value = "  +4.02%  "

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.095"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +3.16%  "

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.591"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +1.38%  "

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "122.56"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +0.91%  "

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9731"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +4.94%  "

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09569"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +1.83%  "

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.449"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("E35").Value = "  -0.17%  "

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.326"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +0.83%  "

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06118"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +2.07%  "

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02263"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +1.61%  "

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.234"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +0.66%  "

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.154"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  -2.53%  "

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.6022"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +1.20%  "

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -0.38%  "

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1901"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +1.86%  "

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.25"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("E45").Value = "  +0.52%  "

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.243"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -2.80%  "

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.20"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  +0.11%  "

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.403"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +1.06%  "

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.935"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -0.18%  "

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00000000318"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +10.54%  "

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06817"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -0.56%  "
